$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.665.59"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "'2.293.55"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'96.40"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "'269.20"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").Value = "'45.33"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "'7.88"
$ws.Range("E12").Value = "  -3.34%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "'15.76"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").Value = "'2.636.33"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "'0.854"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "'2.299.70"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "'43.644.41"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").Value = "'72.00"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'2.53"
$ws.Range("E22").Value = "  +11.06%  "
$ws.Range("D23").Value = "'232.57"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("D24").Value = "'9.09"
$ws.Range("E24").Value = "  -5.78%  "
$ws.Range("D25").Value = "'2.69"
$ws.Range("E25").Value = "  +6.97%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'11.24"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "'3.46"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'38.55"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "'174.83"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "'21.82"
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "'5.41"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("D36").Value = "'4.49"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").Value = "'3.39"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "'12.15"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "'64.79"
$ws.Range("E43").Value = "  +4.28%  "
$ws.Range("D44").Value = "'1.33"
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "'5.16"
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.102"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'97.38"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.443"
$ws.Range("E50").Value = "  +5.00%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.52"
$ws.Range("E51").Value = "  +11.40%  "
